# Shift 5 grouped shapes on slide 3 to the right (x offset only, y unchanged).
# Delta = 173369 EMU (~13.651102 pt) for every group, matching the target diff.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Map of shape Id -> new Left value in points (chosen so that, after the
# host's internal point->EMU conversion, the resulting <a:off x="…"> lands on
# the exact EMU value required by the diff).
$targets = @{
    88  = 467.3346252441406   # 그룹 87  : 5761780 -> 5935149 EMU
    131 = 167.22271728515625  # 그룹 130 : 1950359 -> 2123728 EMU
    134 = 308.5625            # 그룹 133 : 3745374 -> 3918743 EMU
    137 = 467.3202819824219   # 그룹 136 : 5761598 -> 5934967 EMU
    140 = 308.5768127441406   # 그룹 139 : 3745556 -> 3918925 EMU
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($targets.ContainsKey($shape.Id)) {
        $shape.Left = $targets[$shape.Id]
    }
}
